$wb = $excel.ActiveWorkbook

# --- Sheet "Rushing" ---
$ws = $wb.Worksheets.Item("Rushing")

# Row 5: E.Mitchell
$ws.Range("C5").Value = 147
$ws.Range("D5").Value = 77
$ws.Range("E5").Value = 15
$ws.Range("F5").Value = 22

# Row 7: K.Juszczyk
$ws.Range("D7").Value = 3

# Row 9: D.Samuel
$ws.Range("C9").Value = 37
$ws.Range("D9").Value = 33
$ws.Range("E9").Value = 10
$ws.Range("F9").Value = 16

# --- Sheet "Receiving" ---
$ws2 = $wb.Worksheets.Item("Receiving")

# Row 3: E.Mitchell
$ws2.Range("C3").Value = 24
$ws2.Range("D3").Value = 21

# Row 7: D.Samuel
$ws2.Range("C7").Value = 97
$ws2.Range("D7").Value = 59

# Row 8: B.Aiyuk
$ws2.Range("C8").Value = 61

# Row 11: J.Jennings
$ws2.Range("C11").Value = 38
$ws2.Range("D11").Value = 24

# Row 13: G.Kittle
$ws2.Range("C13").Value = 96
$ws2.Range("D13").Value = 77
$ws2.Range("E13").Value = 31
$ws2.Range("F13").Value = 21
$ws2.Range("G13").Value = 9

# Row 15: C.Woerner
$ws2.Range("C15").Value = 6
